# This script rewrites the per-trial stimulus rows (rows 2-42) of Sheet1 so that
# the 1000-subject-version file uses the 'kitchens' memory-task trial order/values
# for subject slot renumbered 41-81 (trial_total = row + 39), re-using the answer
# metadata (cond_cat/cond_mem/correct_answer/stimulus/conceptual/perceptual/
# typicality/n/p_*/r_*) from a permuted set of the 20 unique stimuli, and swaps
# the catch-trial stimulus from row 28 (catch_23.jpg) into row 20 (catch_03.jpg),
# while row 28 itself is populated with the (now unused) img_ensho.png trial data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: trial_total -> 41, stimulus -> 'stimuli/img_c79r7.png'
$ws.Range("F2").Value = 41
$ws.Range("L2").Value = 'stimuli/img_c79r7.png'
$ws.Range("M2").Value = 56.26470588235294
$ws.Range("N2").Value = 34.26470588235294
$ws.Range("O2").Value = 45.26470588235294
$ws.Range("P2").Value = 34
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 2

# Row 3: trial_total -> 42, stimulus -> 'stimuli/img_0j24m.png'
$ws.Range("F3").Value = 42
$ws.Range("L3").Value = 'stimuli/img_0j24m.png'
$ws.Range("M3").Value = 63.6969696969697
$ws.Range("N3").Value = 35.75757575757576
$ws.Range("O3").Value = 49.72727272727273
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 3
$ws.Range("V3").Value = 3

# Row 4: trial_total -> 43, stimulus -> 'stimuli/img_l9t30.png'
$ws.Range("F4").Value = 43
$ws.Range("I4").Value = 'target'
$ws.Range("J4").Value = 'old'
$ws.Range("K4").Value = 'j'
$ws.Range("L4").Value = 'stimuli/img_l9t30.png'
$ws.Range("M4").Value = 67.2
$ws.Range("N4").Value = 43.14285714285715
$ws.Range("O4").Value = 55.17142857142858
$ws.Range("P4").Value = 35
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 4
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 4

# Row 5: trial_total -> 44, stimulus -> 'stimuli/img_fhm45.png'
$ws.Range("F5").Value = 44
$ws.Range("I5").Value = 'target'
$ws.Range("J5").Value = 'old'
$ws.Range("K5").Value = 'j'
$ws.Range("L5").Value = 'stimuli/img_fhm45.png'
$ws.Range("M5").Value = 76.75
$ws.Range("N5").Value = 57.71875
$ws.Range("O5").Value = 67.234375
$ws.Range("P5").Value = 32
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 7
$ws.Range("T5").Value = 7
$ws.Range("U5").Value = 7
$ws.Range("V5").Value = 7

# Row 6: trial_total -> 45, stimulus -> 'stimuli/img_bwo9g.png'
$ws.Range("F6").Value = 45
$ws.Range("J6").Value = 'new'
$ws.Range("K6").Value = 'f'
$ws.Range("L6").Value = 'stimuli/img_bwo9g.png'
$ws.Range("M6").Value = 64.81818181818181
$ws.Range("N6").Value = 42.36363636363637
$ws.Range("O6").Value = 53.59090909090909
$ws.Range("P6").Value = 33
$ws.Range("Q6").Value = 4
$ws.Range("R6").Value = 4
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = 4
$ws.Range("V6").Value = 4
$ws.Range("I6").ClearContents()

# Row 7: trial_total -> 46, stimulus -> 'stimuli/img_01w8b.png'
$ws.Range("F7").Value = 46
$ws.Range("I7").Value = 'target'
$ws.Range("J7").Value = 'old'
$ws.Range("K7").Value = 'j'
$ws.Range("L7").Value = 'stimuli/img_01w8b.png'
$ws.Range("M7").Value = 78.91891891891892
$ws.Range("N7").Value = 61.21621621621622
$ws.Range("O7").Value = 70.06756756756756
$ws.Range("P7").Value = 37
$ws.Range("Q7").Value = 8
$ws.Range("R7").Value = 8
$ws.Range("S7").Value = 8
$ws.Range("T7").Value = 8
$ws.Range("U7").Value = 8
$ws.Range("V7").Value = 8

# Row 8: trial_total -> 47, stimulus -> 'stimuli/img_s9are.png'
$ws.Range("F8").Value = 47
$ws.Range("J8").Value = 'new'
$ws.Range("K8").Value = 'f'
$ws.Range("L8").Value = 'stimuli/img_s9are.png'
$ws.Range("M8").Value = 90.14285714285714
$ws.Range("N8").Value = 75.22857142857143
$ws.Range("O8").Value = 82.68571428571428
$ws.Range("Q8").Value = 10
$ws.Range("R8").Value = 10
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 10
$ws.Range("U8").Value = 10
$ws.Range("V8").Value = 10
$ws.Range("I8").ClearContents()

# Row 9: trial_total -> 48, stimulus -> 'stimuli/img_jz3kd.png'
$ws.Range("F9").Value = 48
$ws.Range("J9").Value = 'new'
$ws.Range("K9").Value = 'f'
$ws.Range("L9").Value = 'stimuli/img_jz3kd.png'
$ws.Range("M9").Value = 72.79411764705883
$ws.Range("N9").Value = 51.64705882352941
$ws.Range("O9").Value = 62.22058823529412
$ws.Range("P9").Value = 34
$ws.Range("Q9").Value = 6
$ws.Range("R9").Value = 6
$ws.Range("S9").Value = 6
$ws.Range("T9").Value = 6
$ws.Range("U9").Value = 6
$ws.Range("V9").Value = 6
$ws.Range("I9").ClearContents()

# Row 10: trial_total -> 49, stimulus -> 'stimuli/img_a8wvq.png'
$ws.Range("F10").Value = 49
$ws.Range("L10").Value = 'stimuli/img_a8wvq.png'
$ws.Range("M10").Value = 86.25925925925925
$ws.Range("N10").Value = 66.25925925925925
$ws.Range("O10").Value = 76.25925925925925
$ws.Range("P10").Value = 27
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = 10
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 10
$ws.Range("U10").Value = 10
$ws.Range("V10").Value = 10

# Row 11: trial_total -> 50, stimulus -> 'stimuli/img_oz18d.png'
$ws.Range("F11").Value = 50
$ws.Range("L11").Value = 'stimuli/img_oz18d.png'
$ws.Range("M11").Value = 78.93939393939394
$ws.Range("N11").Value = 61.03030303030303
$ws.Range("O11").Value = 69.98484848484848
$ws.Range("P11").Value = 33

# Row 12: trial_total -> 51, stimulus -> 'stimuli/img_8fpog.png'
$ws.Range("F12").Value = 51
$ws.Range("L12").Value = 'stimuli/img_8fpog.png'
$ws.Range("M12").Value = 85.41666666666667
$ws.Range("N12").Value = 72.30555555555556
$ws.Range("O12").Value = 78.86111111111111
$ws.Range("P12").Value = 36
$ws.Range("Q12").Value = 10
$ws.Range("R12").Value = 10
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 10
$ws.Range("U12").Value = 10
$ws.Range("V12").Value = 10

# Row 13: trial_total -> 52, stimulus -> 'stimuli/img_ua9bs.png'
$ws.Range("F13").Value = 52

# Row 14: trial_total -> 53, stimulus -> 'stimuli/img_k3abb.png'
$ws.Range("F14").Value = 53
$ws.Range("I14").Value = 'target'
$ws.Range("J14").Value = 'old'
$ws.Range("K14").Value = 'j'
$ws.Range("L14").Value = 'stimuli/img_k3abb.png'
$ws.Range("M14").Value = 35.54054054054054
$ws.Range("N14").Value = 16.54054054054054
$ws.Range("O14").Value = 26.04054054054054
$ws.Range("P14").Value = 37
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 1
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 1
$ws.Range("U14").Value = 1
$ws.Range("V14").Value = 1

# Row 15: trial_total -> 54, stimulus -> 'stimuli/img_uwv6y.png'
$ws.Range("F15").Value = 54
$ws.Range("J15").Value = 'new'
$ws.Range("K15").Value = 'f'
$ws.Range("L15").Value = 'stimuli/img_uwv6y.png'
$ws.Range("M15").Value = 78.88888888888889
$ws.Range("N15").Value = 59.30555555555556
$ws.Range("O15").Value = 69.09722222222223
$ws.Range("P15").Value = 36
$ws.Range("Q15").Value = 8
$ws.Range("R15").Value = 8
$ws.Range("S15").Value = 8
$ws.Range("T15").Value = 8
$ws.Range("U15").Value = 8
$ws.Range("V15").Value = 8
$ws.Range("I15").ClearContents()

# Row 16: trial_total -> 55, stimulus -> 'stimuli/img_84s7n.png'
$ws.Range("F16").Value = 55
$ws.Range("J16").Value = 'new'
$ws.Range("K16").Value = 'f'
$ws.Range("L16").Value = 'stimuli/img_84s7n.png'
$ws.Range("M16").Value = 11.03125
$ws.Range("N16").Value = 2.90625
$ws.Range("O16").Value = 6.96875
$ws.Range("P16").Value = 32
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = 1
$ws.Range("V16").Value = 1
$ws.Range("I16").ClearContents()

# Row 17: trial_total -> 56, stimulus -> 'stimuli/img_i2k07.png'
$ws.Range("F17").Value = 56
$ws.Range("L17").Value = 'stimuli/img_i2k07.png'
$ws.Range("M17").Value = 64.25925925925925
$ws.Range("N17").Value = 40.92592592592592
$ws.Range("O17").Value = 52.59259259259259
$ws.Range("P17").Value = 27
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = 3
$ws.Range("S17").Value = 3
$ws.Range("T17").Value = 3
$ws.Range("U17").Value = 3
$ws.Range("V17").Value = 3

# Row 18: trial_total -> 57, stimulus -> 'stimuli/img_faly8.png'
$ws.Range("F18").Value = 57
$ws.Range("L18").Value = 'stimuli/img_faly8.png'
$ws.Range("M18").Value = 33.41176470588236
$ws.Range("N18").Value = 19.23529411764706
$ws.Range("O18").Value = 26.32352941176471
$ws.Range("P18").Value = 34
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 1
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 1
$ws.Range("U18").Value = 1
$ws.Range("V18").Value = 1

# Row 19: trial_total -> 58, stimulus -> 'stimuli/img_7pgd2.png'
$ws.Range("F19").Value = 58
$ws.Range("J19").Value = 'new'
$ws.Range("K19").Value = 'f'
$ws.Range("L19").Value = 'stimuli/img_7pgd2.png'
$ws.Range("M19").Value = 78.59375
$ws.Range("N19").Value = 57.84375
$ws.Range("O19").Value = 68.21875
$ws.Range("P19").Value = 32
$ws.Range("Q19").Value = 8
$ws.Range("R19").Value = 7
$ws.Range("S19").Value = 7
$ws.Range("T19").Value = 7
$ws.Range("U19").Value = 7
$ws.Range("V19").Value = 7
$ws.Range("I19").ClearContents()

# Row 20: trial_total -> 59, stimulus -> 'stimuli/catch_03.jpg'
$ws.Range("F20").Value = 59
$ws.Range("J20").Value = 'catch'
$ws.Range("L20").Value = 'stimuli/catch_03.jpg'
$ws.Range("H20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("O20").ClearContents()
$ws.Range("P20").ClearContents()
$ws.Range("Q20").ClearContents()
$ws.Range("R20").ClearContents()
$ws.Range("S20").ClearContents()
$ws.Range("T20").ClearContents()
$ws.Range("U20").ClearContents()
$ws.Range("V20").ClearContents()

# Row 21: trial_total -> 60, stimulus -> 'stimuli/img_qbkdt.png'
$ws.Range("F21").Value = 60
$ws.Range("I21").Value = 'target'
$ws.Range("J21").Value = 'old'
$ws.Range("K21").Value = 'j'
$ws.Range("L21").Value = 'stimuli/img_qbkdt.png'
$ws.Range("M21").Value = 69.45714285714286
$ws.Range("N21").Value = 50.02857142857143
$ws.Range("O21").Value = 59.74285714285715
$ws.Range("P21").Value = 35
$ws.Range("Q21").Value = 5
$ws.Range("R21").Value = 5
$ws.Range("S21").Value = 5
$ws.Range("T21").Value = 5
$ws.Range("U21").Value = 5
$ws.Range("V21").Value = 5

# Row 22: trial_total -> 61, stimulus -> 'stimuli/img_xti0z.png'
$ws.Range("F22").Value = 61
$ws.Range("L22").Value = 'stimuli/img_xti0z.png'
$ws.Range("M22").Value = 81.40625
$ws.Range("N22").Value = 61.4375
$ws.Range("O22").Value = 71.421875
$ws.Range("P22").Value = 32
$ws.Range("Q22").Value = 8
$ws.Range("R22").Value = 8
$ws.Range("S22").Value = 8
$ws.Range("T22").Value = 8
$ws.Range("U22").Value = 8
$ws.Range("V22").Value = 8

# Row 23: trial_total -> 62, stimulus -> 'stimuli/img_57os5.png'
$ws.Range("F23").Value = 62
$ws.Range("I23").Value = 'target'
$ws.Range("J23").Value = 'old'
$ws.Range("K23").Value = 'j'
$ws.Range("L23").Value = 'stimuli/img_57os5.png'
$ws.Range("M23").Value = 82.70588235294117
$ws.Range("N23").Value = 65.73529411764706
$ws.Range("O23").Value = 74.22058823529412
$ws.Range("P23").Value = 34
$ws.Range("Q23").Value = 9
$ws.Range("R23").Value = 9
$ws.Range("S23").Value = 9
$ws.Range("T23").Value = 9
$ws.Range("U23").Value = 9
$ws.Range("V23").Value = 9

# Row 24: trial_total -> 63, stimulus -> 'stimuli/img_c0me7.png'
$ws.Range("F24").Value = 63
$ws.Range("L24").Value = 'stimuli/img_c0me7.png'
$ws.Range("M24").Value = 68.4
$ws.Range("N24").Value = 45.62857142857143
$ws.Range("O24").Value = 57.01428571428572
$ws.Range("P24").Value = 35
$ws.Range("Q24").Value = 4
$ws.Range("R24").Value = 4
$ws.Range("S24").Value = 4
$ws.Range("T24").Value = 4
$ws.Range("U24").Value = 4
$ws.Range("V24").Value = 4

# Row 25: trial_total -> 64, stimulus -> 'stimuli/img_qmand.png'
$ws.Range("F25").Value = 64
$ws.Range("I25").Value = 'target'
$ws.Range("J25").Value = 'old'
$ws.Range("K25").Value = 'j'
$ws.Range("L25").Value = 'stimuli/img_qmand.png'
$ws.Range("M25").Value = 86.11764705882354
$ws.Range("N25").Value = 71.02941176470588
$ws.Range("O25").Value = 78.57352941176471
$ws.Range("P25").Value = 34
$ws.Range("Q25").Value = 10
$ws.Range("R25").Value = 10
$ws.Range("S25").Value = 10
$ws.Range("T25").Value = 10
$ws.Range("U25").Value = 10
$ws.Range("V25").Value = 10

# Row 26: trial_total -> 65, stimulus -> 'stimuli/img_5949k.png'
$ws.Range("F26").Value = 65
$ws.Range("I26").Value = 'target'
$ws.Range("J26").Value = 'old'
$ws.Range("K26").Value = 'j'
$ws.Range("L26").Value = 'stimuli/img_5949k.png'
$ws.Range("M26").Value = 60.8
$ws.Range("N26").Value = 39.2
$ws.Range("O26").Value = 50
$ws.Range("P26").Value = 35
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 3
$ws.Range("S26").Value = 3
$ws.Range("T26").Value = 3
$ws.Range("U26").Value = 3
$ws.Range("V26").Value = 3

# Row 27: trial_total -> 66, stimulus -> 'stimuli/img_ifebc.png'
$ws.Range("F27").Value = 66
$ws.Range("L27").Value = 'stimuli/img_ifebc.png'
$ws.Range("M27").Value = 84
$ws.Range("N27").Value = 65.88235294117646
$ws.Range("O27").Value = 74.94117647058823
$ws.Range("P27").Value = 34
$ws.Range("R27").Value = 9
$ws.Range("S27").Value = 9
$ws.Range("T27").Value = 9
$ws.Range("U27").Value = 9
$ws.Range("V27").Value = 9

# Row 28: trial_total -> 67, stimulus -> 'stimuli/img_ensho.png'
$ws.Range("F28").Value = 67
$ws.Range("H28").Value = 'kitchens'
$ws.Range("I28").Value = 'target'
$ws.Range("J28").Value = 'old'
$ws.Range("K28").Value = 'j'
$ws.Range("L28").Value = 'stimuli/img_ensho.png'
$ws.Range("M28").Value = 72.7948717948718
$ws.Range("N28").Value = 54.56410256410256
$ws.Range("O28").Value = 63.67948717948718
$ws.Range("P28").Value = 39
$ws.Range("Q28").Value = 6
$ws.Range("R28").Value = 6
$ws.Range("S28").Value = 6
$ws.Range("T28").Value = 6
$ws.Range("U28").Value = 6
$ws.Range("V28").Value = 6

# Row 29: trial_total -> 68, stimulus -> 'stimuli/img_ncr40.png'
$ws.Range("F29").Value = 68
$ws.Range("J29").Value = 'new'
$ws.Range("K29").Value = 'f'
$ws.Range("L29").Value = 'stimuli/img_ncr40.png'
$ws.Range("M29").Value = 75.66666666666667
$ws.Range("N29").Value = 54.27272727272727
$ws.Range("O29").Value = 64.96969696969697
$ws.Range("P29").Value = 33
$ws.Range("Q29").Value = 6
$ws.Range("R29").Value = 6
$ws.Range("S29").Value = 6
$ws.Range("T29").Value = 6
$ws.Range("U29").Value = 6
$ws.Range("V29").Value = 6
$ws.Range("I29").ClearContents()

# Row 30: trial_total -> 69, stimulus -> 'stimuli/img_05flq.png'
$ws.Range("F30").Value = 69
$ws.Range("L30").Value = 'stimuli/img_05flq.png'
$ws.Range("M30").Value = 47.10344827586207
$ws.Range("N30").Value = 25.72413793103448
$ws.Range("O30").Value = 36.41379310344828
$ws.Range("P30").Value = 29
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = 1
$ws.Range("S30").Value = 1
$ws.Range("T30").Value = 1
$ws.Range("U30").Value = 1
$ws.Range("V30").Value = 1

# Row 31: trial_total -> 70, stimulus -> 'stimuli/img_411xa.png'
$ws.Range("F31").Value = 70
$ws.Range("J31").Value = 'new'
$ws.Range("K31").Value = 'f'
$ws.Range("L31").Value = 'stimuli/img_411xa.png'
$ws.Range("M31").Value = 51.03030303030303
$ws.Range("N31").Value = 28.93939393939394
$ws.Range("O31").Value = 39.98484848484848
$ws.Range("P31").Value = 33
$ws.Range("I31").ClearContents()

# Row 32: trial_total -> 71, stimulus -> 'stimuli/img_uspja.png'
$ws.Range("F32").Value = 71
$ws.Range("L32").Value = 'stimuli/img_uspja.png'
$ws.Range("M32").Value = 54.90909090909091
$ws.Range("N32").Value = 29.12121212121212
$ws.Range("O32").Value = 42.01515151515152
$ws.Range("P32").Value = 33
$ws.Range("Q32").Value = 2
$ws.Range("R32").Value = 2
$ws.Range("S32").Value = 2
$ws.Range("T32").Value = 2
$ws.Range("U32").Value = 2
$ws.Range("V32").Value = 2

# Row 33: trial_total -> 72, stimulus -> 'stimuli/img_j5rpx.png'
$ws.Range("F33").Value = 72
$ws.Range("J33").Value = 'new'
$ws.Range("K33").Value = 'f'
$ws.Range("L33").Value = 'stimuli/img_j5rpx.png'
$ws.Range("M33").Value = 72.24242424242425
$ws.Range("N33").Value = 50
$ws.Range("O33").Value = 61.12121212121212
$ws.Range("P33").Value = 33
$ws.Range("Q33").Value = 5
$ws.Range("R33").Value = 5
$ws.Range("S33").Value = 5
$ws.Range("T33").Value = 5
$ws.Range("U33").Value = 5
$ws.Range("V33").Value = 5
$ws.Range("I33").ClearContents()

# Row 34: trial_total -> 73, stimulus -> 'stimuli/img_as3da.png'
$ws.Range("F34").Value = 73
$ws.Range("L34").Value = 'stimuli/img_as3da.png'
$ws.Range("M34").Value = 84.53125
$ws.Range("N34").Value = 63
$ws.Range("O34").Value = 73.765625
$ws.Range("P34").Value = 32
$ws.Range("Q34").Value = 9
$ws.Range("R34").Value = 9
$ws.Range("S34").Value = 9
$ws.Range("T34").Value = 9
$ws.Range("U34").Value = 9
$ws.Range("V34").Value = 9

# Row 35: trial_total -> 74, stimulus -> 'stimuli/img_cv9qj.png'
$ws.Range("F35").Value = 74
$ws.Range("J35").Value = 'new'
$ws.Range("K35").Value = 'f'
$ws.Range("L35").Value = 'stimuli/img_cv9qj.png'
$ws.Range("M35").Value = 60.34375
$ws.Range("N35").Value = 35.34375
$ws.Range("O35").Value = 47.84375
$ws.Range("Q35").Value = 3
$ws.Range("R35").Value = 3
$ws.Range("S35").Value = 3
$ws.Range("T35").Value = 3
$ws.Range("U35").Value = 3
$ws.Range("V35").Value = 3
$ws.Range("I35").ClearContents()

# Row 36: trial_total -> 75, stimulus -> 'stimuli/img_u9f9l.png'
$ws.Range("F36").Value = 75
$ws.Range("I36").Value = 'target'
$ws.Range("J36").Value = 'old'
$ws.Range("K36").Value = 'j'
$ws.Range("L36").Value = 'stimuli/img_u9f9l.png'
$ws.Range("M36").Value = 77.78571428571429
$ws.Range("N36").Value = 57.25
$ws.Range("O36").Value = 67.51785714285714
$ws.Range("P36").Value = 28
$ws.Range("Q36").Value = 7
$ws.Range("R36").Value = 7
$ws.Range("S36").Value = 7
$ws.Range("T36").Value = 7
$ws.Range("U36").Value = 7
$ws.Range("V36").Value = 7

# Row 37: trial_total -> 76, stimulus -> 'stimuli/img_t1cr9.png'
$ws.Range("F37").Value = 76
$ws.Range("L37").Value = 'stimuli/img_t1cr9.png'
$ws.Range("M37").Value = 73.66666666666667
$ws.Range("N37").Value = 53.51515151515152
$ws.Range("O37").Value = 63.59090909090909
$ws.Range("P37").Value = 33
$ws.Range("Q37").Value = 6
$ws.Range("R37").Value = 6
$ws.Range("S37").Value = 6
$ws.Range("T37").Value = 6
$ws.Range("U37").Value = 6
$ws.Range("V37").Value = 6

# Row 38: trial_total -> 77, stimulus -> 'stimuli/img_z5osu.png'
$ws.Range("F38").Value = 77
$ws.Range("I38").Value = 'target'
$ws.Range("J38").Value = 'old'
$ws.Range("K38").Value = 'j'
$ws.Range("L38").Value = 'stimuli/img_z5osu.png'
$ws.Range("M38").Value = 71.42857142857143
$ws.Range("N38").Value = 47.34285714285714
$ws.Range("O38").Value = 59.38571428571429
$ws.Range("Q38").Value = 5
$ws.Range("R38").Value = 5
$ws.Range("S38").Value = 5
$ws.Range("T38").Value = 5
$ws.Range("U38").Value = 5
$ws.Range("V38").Value = 5

# Row 39: trial_total -> 78, stimulus -> 'stimuli/img_uy1n4.png'
$ws.Range("F39").Value = 78
$ws.Range("J39").Value = 'new'
$ws.Range("K39").Value = 'f'
$ws.Range("L39").Value = 'stimuli/img_uy1n4.png'
$ws.Range("M39").Value = 76.30555555555556
$ws.Range("N39").Value = 55.33333333333334
$ws.Range("O39").Value = 65.81944444444444
$ws.Range("P39").Value = 36
$ws.Range("Q39").Value = 7
$ws.Range("R39").Value = 7
$ws.Range("S39").Value = 7
$ws.Range("T39").Value = 7
$ws.Range("U39").Value = 7
$ws.Range("V39").Value = 7
$ws.Range("I39").ClearContents()

# Row 40: trial_total -> 79, stimulus -> 'stimuli/img_h1yyu.png'
$ws.Range("F40").Value = 79
$ws.Range("I40").Value = 'target'
$ws.Range("J40").Value = 'old'
$ws.Range("K40").Value = 'j'
$ws.Range("L40").Value = 'stimuli/img_h1yyu.png'
$ws.Range("M40").Value = 64.8529411764706
$ws.Range("N40").Value = 46.61764705882353
$ws.Range("O40").Value = 55.73529411764706
$ws.Range("P40").Value = 34
$ws.Range("Q40").Value = 4
$ws.Range("R40").Value = 4
$ws.Range("S40").Value = 4
$ws.Range("T40").Value = 4
$ws.Range("U40").Value = 4
$ws.Range("V40").Value = 4

# Row 41: trial_total -> 80, stimulus -> 'stimuli/img_hfz8w.png'
$ws.Range("F41").Value = 80
$ws.Range("L41").Value = 'stimuli/img_hfz8w.png'
$ws.Range("M41").Value = 55.46153846153846
$ws.Range("N41").Value = 27.28205128205128
$ws.Range("O41").Value = 41.37179487179487
$ws.Range("P41").Value = 39
$ws.Range("Q41").Value = 2
$ws.Range("R41").Value = 2
$ws.Range("S41").Value = 2
$ws.Range("T41").Value = 2
$ws.Range("U41").Value = 2
$ws.Range("V41").Value = 2

# Row 42: trial_total -> 81, stimulus -> 'stimuli/img_xesl0.png'
$ws.Range("F42").Value = 81
$ws.Range("J42").Value = 'new'
$ws.Range("K42").Value = 'f'
$ws.Range("L42").Value = 'stimuli/img_xesl0.png'
$ws.Range("M42").Value = 69.28571428571429
$ws.Range("N42").Value = 47.35714285714285
$ws.Range("O42").Value = 58.32142857142857
$ws.Range("P42").Value = 28
$ws.Range("Q42").Value = 5
$ws.Range("R42").Value = 5
$ws.Range("S42").Value = 5
$ws.Range("T42").Value = 5
$ws.Range("U42").Value = 5
$ws.Range("V42").Value = 5
$ws.Range("I42").ClearContents()

